$d = $word.ActiveDocument

$map = @{
    "127904385" = "127905004"
    "127904386" = "127905005"
    "127904387" = "127905006"
    "127904388" = "127905007"
    "127904389" = "127905008"
    "127904390" = "127905009"
}

foreach ($old in $map.Keys) {
    $new = $map[$old]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
